# [new feature - search product by name or id] added to the SIV
#
# This script updates the "Sheet" (product inventory) and "Logs" (stock
# change history) worksheets to reflect:
#   - Product #1 "doll" (was Iphone8 / id 2) with stock 5, last updated 2025-04-20
#   - Product #2 "Pen toy" (was IphoneX / id 1) with stock 10, status OK
#   - The Logs sheet trimmed down to a single remaining entry (row 2), whose
#     date was bumped to 2025-04-20 and annotated with the user "hengty"
#   - Older log rows (3-9) removed, since they no longer apply

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Sheet" worksheet - product list
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet")

# Row 2: id 1, doll, stock 5, last updated 2025-04-20
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "1"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "doll"
$ws1.Range("C2").Value = 5
$ws1.Range("E2").NumberFormat = "@"
$ws1.Range("E2").Value = "2025-04-20"

# Row 3: id 2, Pen toy, stock 10, status OK
$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "2"
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "Pen toy"
$ws1.Range("C3").Value = 10
$ws1.Range("F3").NumberFormat = "@"
$ws1.Range("F3").Value = "OK"

$ws1.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# "Logs" worksheet - stock change history
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Logs")

# Update the remaining log entry's date
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "2025-04-20"

# Record who made the change
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "hengty"

# Remove the now-obsolete log rows (3 through 9)
$ws2.Range("A3:A9").EntireRow.Delete()

$ws2.Range("A1").Select() | Out-Null

$ws1.Activate() | Out-Null
